# Konfliktuskezelés — add the "halálos karambol" story paragraph text
# and justify the paragraph it lives in.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Justify the body paragraph (adds <w:jc w:val="both"/> to its pPr).
#    It's the 3rd paragraph in the document (heading, blank line, body).
# ---------------------------------------------------------------------
$d.Paragraphs(3).Alignment = 3   # wdAlignParagraphJustify -> w:jc val="both"

# ---------------------------------------------------------------------
# 2) Append the new story text right after "...a közelmúltból", before
#    the existing _GoBack bookmark, keeping the surrounding Times New
#    Roman / 12pt run formatting.
# ---------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("közelmúltból", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)

$rng.Text = ". 2020 őszén mentünk a rokonokhoz édesapámmal."
$rng.Collapse(0)

$rng.Text = " A "
$rng.Collapse(0)

$rng.Text = "Bajáról kivezető úton haladtunk"
$rng.Collapse(0)

$rng.Text = ","
$rng.Collapse(0)

$rng.Text = " amikor egy kisebb dugó miatt meg kellett állnunk. Mint kiderült, egy halálos karambol "
$rng.Collapse(0)

$rng.Text = "miatt alakult ki torlódás, amelyben ketten meghaltak. A kocsisorban egy úr a kocsijával úgy gondolta, hogy"
$rng.Collapse(0)

$rng.Text = " hátra "
$rng.Collapse(0)

$rng.Text = "tolat"
$rng.Collapse(0)

$rng.Text = ","
$rng.Collapse(0)

$rng.Text = " mert az vicces. Természetesen nekünk is koccant, és mivel az a nap amúgy sem volt az év napja, kicsit felforrt az agyvizünk (inkább édesapámé). Ki is szállt és elkezdett verbális"
$rng.Collapse(0)

$rng.Text = "an "
$rng.Collapse(0)

$rng.Text = "kommunikálni a férfival, aki inkább az autójában maradt. "
$rng.Collapse(0)

$rng.Text = "Miután a kötelező papírokat kitöltötték, megegyezve a következményekben,"
$rng.Collapse(0)

$rng.Text = " indultunk tovább. Kikerülve a dugóból"
$rng.Collapse(0)

$rng.Text = ", "
$rng.Collapse(0)

$rng.Text = "kb. fél óra múlva "
$rng.Collapse(0)

$rng.Text = "kaptam a hírt, hogy az egyik osztálytársam , aki közel ült az órákon, "
$rng.Collapse(0)

$rng.Text = "COVIDos lett, így nem tudtunk találkozni a családtagokkal. Hazafelé elmentünk a mi koccanásunkban vétkes sofőrhöz és elrendeztük a maradék papírt és kártérítést"
$rng.Collapse(0)
